# Re-brand the Insight StructureDefinition workbook from "ibm.com" /
# "Alvearie Team" to "linuxforhealth.org" / "LinuxForHealth Team", bump the
# version and publication date (LinuxForHealth/alvearie-fhir-ig rebrand).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Metadata" sheet - simple Property / Value table
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------------------
# "Elements" sheet - FHIR StructureDefinition element grid. Only the
# "ibm.com" URLs embedded in the Type(s)/Fixed Value/Binding Value Set
# columns need the domain swapped; row/column layout is unchanged.
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.extension slices (Type(s) column J) - each is "Extension {url}\n"
$elements.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/path}`n"
$elements.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-id}`n"
$elements.Range("J7").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/source-id}`n"
$elements.Range("J8").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/detected}`n"
$elements.Range("J9").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/category}`n"
$elements.Range("J15").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-type}`n"
$elements.Range("J16").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-detail}`n"
$elements.Range("J17").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-result-summary}`n"

# Fixed Value column (Q) for the url-typed children
$elements.Range("Q12").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/category"
$elements.Range("Q18").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight"

# Binding Value Set column (Y)
$elements.Range("Y14").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/insight-category-values"

# The Binding Value Set column widens to fit the longer linuxforhealth.org
# URL - let Excel recompute the best-fit width for column Y.
$elements.Columns.Item(25).AutoFit()
